$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text (matching original inlineStr formatting)
$textCells = @("D5", "D6", "D7", "D11", "D12", "D15", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated price / volume values
$ws.Range('D2').Value = '67.002.36'
$ws.Range('E2').Value = '  +4.38%  '
$ws.Range('D3').Value = '3.265.03'
$ws.Range('E3').Value = '  +2.81%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '579.48'
$ws.Range('E5').Value = '  +3.01%  '
$ws.Range('D6').Value = '177.45'
$ws.Range('E6').Value = '  +4.10%  '
$ws.Range('D7').Value = '0.607'
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '3.264.22'
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').Value = '6.73'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  +4.32%  '
$ws.Range('D13').Value = '3.833.55'
$ws.Range('E13').Value = '  +2.91%  '
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '28.19'
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').Value = '67.015.61'
$ws.Range('E16').Value = '  +4.40%  '
$ws.Range('D17').Value = '0.0000168'
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('D18').Value = '3.266.37'
$ws.Range('E18').Value = '  +3.01%  '
$ws.Range('D19').Value = '5.85'
$ws.Range('E19').Value = '  +3.17%  '
$ws.Range('D20').Value = '13.45'
$ws.Range('E20').Value = '  +2.81%  '
$ws.Range('D21').Value = '370.63'
$ws.Range('E21').Value = '  +5.21%  '
$ws.Range('D22').Value = '7.66'
$ws.Range('E22').Value = '  +6.48%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '71.07'
$ws.Range('E24').Value = '  +3.31%  '
$ws.Range('D25').Value = '0.511'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').Value = '3.404.73'
$ws.Range('E26').Value = '  +2.63%  '
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('D28').Value = '9.83'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('E29').Value = '  +2.47%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '1.99'
$ws.Range('E31').Value = '  +5.11%  '
$ws.Range('D32').Value = '5.65'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').Value = '22.62'
$ws.Range('E33').Value = '  +2.42%  '
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  +4.32%  '
$ws.Range('D36').Value = '6.82'
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('D37').Value = '168.28'
$ws.Range('E37').Value = '  +8.46%  '
$ws.Range('D38').Value = '1.51'
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('D39').Value = '0.859'
$ws.Range('E39').Value = '  +6.54%  '
$ws.Range('D40').Value = '1.86'
$ws.Range('E40').Value = '  +10.14%  '
$ws.Range('D41').Value = '27.37'
$ws.Range('E41').Value = '  +5.83%  '
$ws.Range('D42').Value = '2.767.36'
$ws.Range('E42').Value = '  +4.52%  '
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('E44').Value = '  +6.89%  '
$ws.Range('D45').Value = '4.37'
$ws.Range('E45').Value = '  +5.01%  '
$ws.Range('D46').Value = '344.42'
$ws.Range('E46').Value = '  +4.48%  '
$ws.Range('D47').Value = '40.49'
$ws.Range('E47').Value = '  +5.03%  '
$ws.Range('D48').Value = '0.0675'
$ws.Range('E48').Value = '  +3.32%  '
$ws.Range('D49').Value = '24.90'
$ws.Range('E49').Value = '  +5.33%  '
$ws.Range('D50').Value = '0.0280'
$ws.Range('E50').Value = '  +3.53%  '
$ws.Range('E51').Value = '  +3.66%  '
